$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 16.52233964717913
$ws.Range("C2").Value = 10.32468481425973
$ws.Range("D2").Value = 11.24921220768562
$ws.Range("F2").Value = 30.52285013937888
$ws.Range("G2").Value = 3.632351806462029
$ws.Range("J2").Value = 11.31980928077249
$ws.Range("O2").Value = 21.96016188730161
$ws.Range("B3").Value = 15.80019147163543
$ws.Range("C3").Value = 9.67976870335449
$ws.Range("D3").Value = 11.14507077160001
$ws.Range("F3").Value = 30.59970734668889
$ws.Range("G3").Value = 3.63477847272315
$ws.Range("J3").Value = 11.30181205426907
$ws.Range("O3").Value = 22.08706669769633
$ws.Range("B4").Value = 15.34005903821558
$ws.Range("C4").Value = 9.260058416481153
$ws.Range("D4").Value = 11.08233381936053
$ws.Range("F4").Value = 30.65804030734518
$ws.Range("G4").Value = 3.636345691123417
$ws.Range("J4").Value = 11.29310745820674
$ws.Range("O4").Value = 22.17252561077112
$ws.Range("B5").Value = 15.14858140521858
$ws.Range("C5").Value = 9.083088502094485
$ws.Range("D5").Value = 11.05709402216911
$ws.Range("F5").Value = 30.68459741780283
$ws.Range("G5").Value = 3.637003832085298
$ws.Range("J5").Value = 11.29015228783193
$ws.Range("O5").Value = 22.20923569260494
$ws.Range("B6").Value = 15.11655449873704
$ws.Range("C6").Value = 9.053345518279089
$ws.Range("D6").Value = 11.05292331497175
$ws.Range("F6").Value = 30.68917498051149
$ws.Range("G6").Value = 3.637114294759297
$ws.Range("J6").Value = 11.28969739481252
$ws.Range("O6").Value = 22.21544489880037
$ws.Range("B7").Value = 15.33749244036927
$ws.Range("C7").Value = 9.257695716303981
$ws.Range("D7").Value = 11.08199207706024
$ws.Range("F7").Value = 30.65838720764662
$ws.Range("G7").Value = 3.636354488055368
$ws.Range("J7").Value = 11.29306520415928
$ws.Range("O7").Value = 22.17301307835808
$ws.Range("B8").Value = 16.27694760059106
$ws.Range("C8").Value = 10.10725880828519
$ws.Range("D8").Value = 11.21306765095578
$ws.Range("F8").Value = 30.54702877652385
$ws.Range("G8").Value = 3.633172530122982
$ws.Range("J8").Value = 11.31311840184809
$ws.Range("O8").Value = 22.00234748418717
$ws.Range("B9").Value = 17.97763938485242
$ws.Range("C9").Value = 11.58405821616121
$ws.Range("D9").Value = 11.47865740996335
$ws.Range("F9").Value = 30.41771040200414
$ws.Range("G9").Value = 3.627542565576452
$ws.Range("J9").Value = 11.37094144348688
$ws.Range("O9").Value = 21.72797180334648
$ws.Range("B10").Value = 19.13098058073487
$ws.Range("C10").Value = 12.55306429647234
$ws.Range("D10").Value = 11.67760702754314
$ws.Range("F10").Value = 30.37777591062527
$ws.Range("G10").Value = 3.623773822278834
$ws.Range("J10").Value = 11.42452127869939
$ws.Range("O10").Value = 21.56380293553112
$ws.Range("B11").Value = 19.63319515857744
$ws.Range("C11").Value = 12.96867955290019
$ws.Range("D11").Value = 11.76864377228149
$ws.Range("F11").Value = 30.37170227292016
$ws.Range("G11").Value = 3.622138249645897
$ws.Range("J11").Value = 11.4512579114088
$ws.Range("O11").Value = 21.49738150681094
$ws.Range("B12").Value = 19.82003745840061
$ws.Range("C12").Value = 13.12244229500781
$ws.Range("D12").Value = 11.80316848569374
$ws.Range("F12").Value = 30.37114950296481
$ws.Range("G12").Value = 3.621530171187384
$ws.Range("J12").Value = 11.4617173464949
$ws.Range("O12").Value = 21.4734285363916
$ws.Range("B13").Value = 19.77994750248023
$ws.Range("C13").Value = 13.08948777207588
$ws.Range("D13").Value = 11.795731105952
$ws.Range("F13").Value = 30.37119073341061
$ws.Range("G13").Value = 3.621660631161038
$ws.Range("J13").Value = 11.45944991050603
$ws.Range("O13").Value = 21.47853371052271
$ws.Range("B14").Value = 19.6486341897678
$ws.Range("C14").Value = 12.98140227655629
$ws.Range("D14").Value = 11.77148326394457
$ws.Range("F14").Value = 30.37162173924854
$ws.Range("G14").Value = 3.62208799699629
$ws.Range("J14").Value = 11.45211173428005
$ws.Range("O14").Value = 21.49538677796692
$ws.Range("B15").Value = 19.56776372031597
$ws.Range("C15").Value = 12.91472520494614
$ws.Range("D15").Value = 11.75663665348463
$ws.Range("F15").Value = 30.37211349157937
$ws.Range("G15").Value = 3.622351237762153
$ws.Range("J15").Value = 11.44766034333685
$ws.Range("O15").Value = 21.50586630633539
$ws.Range("B16").Value = 19.09769838723893
$ws.Range("C16").Value = 12.52539594491199
$ws.Range("D16").Value = 11.67166610258714
$ws.Range("F16").Value = 30.37841692234452
$ws.Range("G16").Value = 3.623882292206973
$ws.Range("J16").Value = 11.42282111315247
$ws.Range("O16").Value = 21.56831093807503
$ws.Range("B17").Value = 18.80349715113559
$ws.Range("C17").Value = 12.2801025654173
$ws.Range("D17").Value = 11.61965778549611
$ws.Range("F17").Value = 30.38538757603535
$ws.Range("G17").Value = 3.624841695906426
$ws.Range("J17").Value = 11.40818487864261
$ws.Range("O17").Value = 21.60874255155012
$ws.Range("B18").Value = 18.63217286722926
$ws.Range("C18").Value = 12.13664498526302
$ws.Range("D18").Value = 11.58979572817546
$ws.Range("F18").Value = 30.3905347597325
$ws.Range("G18").Value = 3.625400944656739
$ws.Range("J18").Value = 11.39998923430985
$ws.Range("O18").Value = 21.63277451074422
$ws.Range("B19").Value = 18.57380701748851
$ws.Range("C19").Value = 12.08766547227382
$ws.Range("D19").Value = 11.57969462094951
$ws.Range("F19").Value = 30.39247265266878
$ws.Range("G19").Value = 3.625591573804565
$ws.Range("J19").Value = 11.39725272708141
$ws.Range("O19").Value = 21.64104437272469
$ws.Range("B20").Value = 18.83503433126049
$ws.Range("C20").Value = 12.30645989222051
$ws.Range("D20").Value = 11.62518898991305
$ws.Range("F20").Value = 30.38452770932035
$ws.Range("G20").Value = 3.624738797713922
$ws.Range("J20").Value = 11.40971991330997
$ws.Range("O20").Value = 21.60435806031028
$ws.Range("B21").Value = 19.68729538839813
$ws.Range("C21").Value = 13.01324787903703
$ws.Range("D21").Value = 11.77860426035357
$ws.Range("F21").Value = 30.37144766758242
$ws.Range("G21").Value = 3.621962163708504
$ws.Range("J21").Value = 11.45425808817818
$ws.Range("O21").Value = 21.49040398060825
$ws.Range("B22").Value = 20.22481750655732
$ws.Range("C22").Value = 13.45406680242445
$ws.Range("D22").Value = 11.87915595555401
$ws.Range("F22").Value = 30.3730855561293
$ws.Range("G22").Value = 3.620213178252928
$ws.Range("J22").Value = 11.48531554384533
$ws.Range("O22").Value = 21.42292461740752
$ws.Range("B23").Value = 19.93974542029292
$ws.Range("C23").Value = 13.22072383515241
$ws.Range("D23").Value = 11.8254718830772
$ws.Range("F23").Value = 30.37127707230546
$ws.Range("G23").Value = 3.621140652646739
$ws.Range("J23").Value = 11.46856298598613
$ws.Range("O23").Value = 21.45829566369592
$ws.Range("B24").Value = 18.82078316668695
$ws.Range("C24").Value = 12.29455132983632
$ws.Range("D24").Value = 11.62268821141127
$ws.Range("F24").Value = 30.38491290563839
$ws.Range("G24").Value = 3.624785294075508
$ws.Range("J24").Value = 11.40902524185283
$ws.Range("O24").Value = 21.60633783723695
$ws.Range("B25").Value = 17.53383044556278
$ws.Range("C25").Value = 11.20489182091859
$ws.Range("D25").Value = 11.40603557097633
$ws.Range("F25").Value = 30.44307315907495
$ws.Range("G25").Value = 3.629000766180936
$ws.Range("J25").Value = 11.35333509318817
$ws.Range("O25").Value = 21.79567437103117
